# Auto-generated Excel COM-interop script applying numeric corrections
# to the Exodus_Profits leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 2296.5386
$ws.Cells.Item(28, 9).Value = 881.75
$ws.Cells.Item(28, 11).Value = 881.75
$ws.Cells.Item(28, 13).Value = -396.75
$ws.Cells.Item(33, 8).Value = 161.08
$ws.Cells.Item(33, 9).Value = 141.66667
$ws.Cells.Item(33, 11).Value = 141.66667
$ws.Cells.Item(33, 13).Value = 87.33332999999999
$ws.Cells.Item(40, 8).Value = 4881.727
$ws.Cells.Item(40, 9).Value = 5060.2
$ws.Cells.Item(40, 11).Value = 5060.2
$ws.Cells.Item(40, 13).Value = -4885.2
$ws.Cells.Item(43, 8).Value = 2424.1
$ws.Cells.Item(43, 9).Value = 2448.5
$ws.Cells.Item(43, 10).Value = 2387.5
$ws.Cells.Item(43, 11).Value = 2448.5
$ws.Cells.Item(43, 12).Value = 2387.5
$ws.Cells.Item(43, 13).Value = -2379.5
$ws.Cells.Item(43, 14).Value = -2525.5
$ws.Cells.Item(98, 8).Value = 484.3846
$ws.Cells.Item(98, 9).Value = 518.2727
$ws.Cells.Item(98, 10).Value = 298
$ws.Cells.Item(98, 11).Value = 518.2727
$ws.Cells.Item(98, 12).Value = 298
$ws.Cells.Item(98, 13).Value = 979.7273
$ws.Cells.Item(98, 14).Value = -3294
$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 14).Value = $null
$ws.Cells.Item(122, 8).Value = 484.3846
$ws.Cells.Item(122, 9).Value = 518.2727
$ws.Cells.Item(122, 10).Value = 298
$ws.Cells.Item(122, 11).Value = 1554.8181
$ws.Cells.Item(122, 12).Value = 894
$ws.Cells.Item(122, 13).Value = 895.1819
$ws.Cells.Item(122, 14).Value = -5794
$ws.Cells.Item(132, 8).Value = 1458.8823
$ws.Cells.Item(132, 9).Value = 1286.7333
$ws.Cells.Item(132, 11).Value = 3860.199900000001
$ws.Cells.Item(132, 13).Value = -1330.199900000001
$ws.Cells.Item(135, 8).Value = 1434.5217
$ws.Cells.Item(135, 9).Value = 1434.5217
$ws.Cells.Item(135, 11).Value = 12910.6953
$ws.Cells.Item(135, 13).Value = -10375.6953
$ws.Cells.Item(137, 8).Value = 855462.5
$ws.Cells.Item(137, 9).Value = 1906.3334
$ws.Cells.Item(137, 10).Value = 1321038.5
$ws.Cells.Item(137, 11).Value = 5719.0002
$ws.Cells.Item(137, 12).Value = 3963115.5
$ws.Cells.Item(137, 13).Value = -3169.0002
$ws.Cells.Item(137, 14).Value = -3968215.5
$ws.Cells.Item(138, 8).Value = 2303.4736
$ws.Cells.Item(138, 10).Value = 4291.3335
$ws.Cells.Item(138, 12).Value = 12874.0005
$ws.Cells.Item(138, 14).Value = -23154.0005
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2537.0588
$ws.Cells.Item(32, 9).Value = 2098.7673
$ws.Cells.Item(32, 11).Value = 2098.7673
$ws.Cells.Item(32, 13).Value = -1811.7673
$ws.Cells.Item(61, 8).Value = 73425.42999999999
$ws.Cells.Item(61, 9).Value = 2096.3
$ws.Cells.Item(61, 11).Value = 2096.3
$ws.Cells.Item(61, 13).Value = -1884.3
$ws.Cells.Item(74, 8).Value = 4421.24
$ws.Cells.Item(74, 9).Value = 1433.1875
$ws.Cells.Item(74, 10).Value = 9733.333000000001
$ws.Cells.Item(74, 11).Value = 1433.1875
$ws.Cells.Item(74, 12).Value = 9733.333000000001
$ws.Cells.Item(74, 13).Value = -559.1875
$ws.Cells.Item(74, 14).Value = -11481.333
$ws.Cells.Item(77, 8).Value = 4421.24
$ws.Cells.Item(77, 9).Value = 1433.1875
$ws.Cells.Item(77, 10).Value = 9733.333000000001
$ws.Cells.Item(77, 11).Value = 7165.9375
$ws.Cells.Item(77, 12).Value = 48666.665
$ws.Cells.Item(77, 13).Value = -2797.9375
$ws.Cells.Item(77, 14).Value = -57402.665
$ws.Cells.Item(97, 8).Value = 1017.3333
$ws.Cells.Item(97, 9).Value = 1017.3333
$ws.Cells.Item(97, 11).Value = 1017.3333
$ws.Cells.Item(97, 13).Value = -521.3333
$ws.Cells.Item(122, 8).Value = 2021.7646
$ws.Cells.Item(122, 10).Value = 1836.75
$ws.Cells.Item(122, 12).Value = 5510.25
$ws.Cells.Item(122, 14).Value = -10410.25
$ws.Cells.Item(132, 8).Value = 2235.1072
$ws.Cells.Item(132, 9).Value = 2128.42
$ws.Cells.Item(132, 10).Value = 3124.1667
$ws.Cells.Item(132, 11).Value = 6385.26
$ws.Cells.Item(132, 12).Value = 9372.500100000001
$ws.Cells.Item(132, 13).Value = -3855.26
$ws.Cells.Item(132, 14).Value = -14432.5001
$ws.Cells.Item(136, 8).Value = 73425.42999999999
$ws.Cells.Item(136, 9).Value = 2096.3
$ws.Cells.Item(136, 11).Value = 6288.900000000001
$ws.Cells.Item(136, 13).Value = -3738.900000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 4588
$ws.Cells.Item(86, 10).Value = 4470.375
$ws.Cells.Item(86, 12).Value = 4470.375
$ws.Cells.Item(86, 14).Value = -6716.375
$ws.Cells.Item(89, 8).Value = 4588
$ws.Cells.Item(89, 10).Value = 4470.375
$ws.Cells.Item(89, 12).Value = 22351.875
$ws.Cells.Item(89, 14).Value = -33583.875
$ws.Cells.Item(99, 8).Value = 113796.336
$ws.Cells.Item(99, 9).Value = 252249.5
$ws.Cells.Item(99, 11).Value = 252249.5
$ws.Cells.Item(99, 13).Value = -250751.5
$ws.Cells.Item(105, 8).Value = 85613.836
$ws.Cells.Item(105, 9).Value = 126639.625
$ws.Cells.Item(105, 11).Value = 126639.625
$ws.Cells.Item(105, 13).Value = -124892.625
$ws.Cells.Item(107, 8).Value = 3750.6365
$ws.Cells.Item(107, 9).Value = 3825.7
$ws.Cells.Item(107, 11).Value = 3825.7
$ws.Cells.Item(107, 13).Value = -1905.7
$ws.Cells.Item(134, 8).Value = 2909.9211
$ws.Cells.Item(134, 9).Value = 1937.5161
$ws.Cells.Item(134, 11).Value = 5812.5483
$ws.Cells.Item(134, 13).Value = -3277.5483
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3000
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 3000
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 3000
$ws.Cells.Item(16, 13).Value = $null
$ws.Cells.Item(16, 14).Value = -3574
$ws.Cells.Item(18, 8).Value = 24991
$ws.Cells.Item(18, 10).Value = 24991
$ws.Cells.Item(18, 12).Value = 24991
$ws.Cells.Item(18, 14).Value = -25451
$ws.Cells.Item(31, 8).Value = 2649.15
$ws.Cells.Item(31, 9).Value = 2261
$ws.Cells.Item(31, 10).Value = 2907.9167
$ws.Cells.Item(31, 11).Value = 2261
$ws.Cells.Item(31, 12).Value = 2907.9167
$ws.Cells.Item(31, 13).Value = -1966
$ws.Cells.Item(31, 14).Value = -3497.9167
$ws.Cells.Item(34, 8).Value = 2649.15
$ws.Cells.Item(34, 9).Value = 2261
$ws.Cells.Item(34, 10).Value = 2907.9167
$ws.Cells.Item(34, 11).Value = 2261
$ws.Cells.Item(34, 12).Value = 2907.9167
$ws.Cells.Item(34, 13).Value = -2059
$ws.Cells.Item(34, 14).Value = -3311.9167
$ws.Cells.Item(58, 8).Value = 2711.8
$ws.Cells.Item(58, 9).Value = 2711.8
$ws.Cells.Item(58, 11).Value = 2711.8
$ws.Cells.Item(58, 13).Value = -2508.8
$ws.Cells.Item(86, 8).Value = 6865.3
$ws.Cells.Item(86, 9).Value = 5074.75
$ws.Cells.Item(86, 10).Value = 8059
$ws.Cells.Item(86, 11).Value = 5074.75
$ws.Cells.Item(86, 12).Value = 8059
$ws.Cells.Item(86, 13).Value = -3951.75
$ws.Cells.Item(86, 14).Value = -10305
$ws.Cells.Item(89, 8).Value = 6865.3
$ws.Cells.Item(89, 9).Value = 5074.75
$ws.Cells.Item(89, 10).Value = 8059
$ws.Cells.Item(89, 11).Value = 25373.75
$ws.Cells.Item(89, 12).Value = 40295
$ws.Cells.Item(89, 13).Value = -19757.75
$ws.Cells.Item(89, 14).Value = -51527
$ws.Cells.Item(96, 8).Value = 21847.5
$ws.Cells.Item(96, 10).Value = 21847.5
$ws.Cells.Item(96, 12).Value = 21847.5
$ws.Cells.Item(96, 14).Value = -27339.5
$ws.Cells.Item(104, 8).Value = 0
$ws.Cells.Item(104, 10).Value = 0
$ws.Cells.Item(104, 12).Value = 0
$ws.Cells.Item(104, 14).Value = $null
$ws.Cells.Item(105, 8).Value = 4204.4
$ws.Cells.Item(105, 9).Value = 2755.5
$ws.Cells.Item(105, 11).Value = 2755.5
$ws.Cells.Item(105, 13).Value = -1008.5
$ws.Cells.Item(109, 8).Value = 43857.145
$ws.Cells.Item(109, 10).Value = 43857.145
$ws.Cells.Item(109, 12).Value = 43857.145
$ws.Cells.Item(109, 14).Value = -45937.145
$ws.Cells.Item(113, 8).Value = 3000
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 3000
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 3000
$ws.Cells.Item(113, 13).Value = $null
$ws.Cells.Item(113, 14).Value = -7340
$ws.Cells.Item(132, 8).Value = 4334738
$ws.Cells.Item(132, 9).Value = 1000
$ws.Cells.Item(132, 11).Value = 3000
$ws.Cells.Item(132, 13).Value = -470
$ws.Cells.Item(134, 8).Value = 31446.4
$ws.Cells.Item(134, 9).Value = 3037.4482
$ws.Cells.Item(134, 11).Value = 9112.3446
$ws.Cells.Item(134, 13).Value = -6577.3446
$ws.Cells.Item(136, 8).Value = 2711.8
$ws.Cells.Item(136, 9).Value = 2711.8
$ws.Cells.Item(136, 11).Value = 8135.400000000001
$ws.Cells.Item(136, 13).Value = -5585.400000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 792.4167
$ws.Cells.Item(5, 9).Value = 567.55554
$ws.Cells.Item(5, 10).Value = 1467
$ws.Cells.Item(5, 11).Value = 1702.66662
$ws.Cells.Item(5, 12).Value = 4401
$ws.Cells.Item(5, 13).Value = -1590.66662
$ws.Cells.Item(5, 14).Value = -4625
$ws.Cells.Item(26, 8).Value = 1530.25
$ws.Cells.Item(26, 9).Value = 373
$ws.Cells.Item(26, 10).Value = 5002
$ws.Cells.Item(26, 11).Value = 1119
$ws.Cells.Item(26, 12).Value = 15006
$ws.Cells.Item(26, 13).Value = -831
$ws.Cells.Item(26, 14).Value = -15582
$ws.Cells.Item(34, 8).Value = 2080.6667
$ws.Cells.Item(34, 10).Value = 3586
$ws.Cells.Item(34, 12).Value = 10758
$ws.Cells.Item(34, 14).Value = -10926
$ws.Cells.Item(39, 8).Value = 5985.0557
$ws.Cells.Item(39, 9).Value = 1199.5
$ws.Cells.Item(39, 10).Value = 6583.25
$ws.Cells.Item(39, 11).Value = 3598.5
$ws.Cells.Item(39, 12).Value = 19749.75
$ws.Cells.Item(39, 13).Value = -3304.5
$ws.Cells.Item(39, 14).Value = -20337.75
$ws.Cells.Item(55, 8).Value = 4899.8
$ws.Cells.Item(55, 10).Value = 9399.6
$ws.Cells.Item(55, 12).Value = 28198.8
$ws.Cells.Item(55, 14).Value = -28552.8
$ws.Cells.Item(61, 8).Value = 2316.1667
$ws.Cells.Item(61, 10).Value = 1873.25
$ws.Cells.Item(61, 12).Value = 5619.75
$ws.Cells.Item(61, 14).Value = -6049.75
$ws.Cells.Item(68, 8).Value = 1698.75
$ws.Cells.Item(68, 10).Value = 1665
$ws.Cells.Item(68, 12).Value = 4995
$ws.Cells.Item(68, 14).Value = -6617
$ws.Cells.Item(70, 8).Value = 12
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 14).Value = $null
$ws.Cells.Item(71, 8).Value = 1698.75
$ws.Cells.Item(71, 10).Value = 1665
$ws.Cells.Item(71, 12).Value = 14985
$ws.Cells.Item(71, 14).Value = -23097
$ws.Cells.Item(73, 8).Value = 12
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 14).Value = $null
$ws.Cells.Item(107, 8).Value = 1959.3572
$ws.Cells.Item(107, 10).Value = 1993.4445
$ws.Cells.Item(107, 12).Value = 5980.333500000001
$ws.Cells.Item(107, 14).Value = -9820.333500000001
$ws.Cells.Item(132, 8).Value = 6528.5
$ws.Cells.Item(132, 10).Value = 6883.1177
$ws.Cells.Item(132, 12).Value = 61948.0593
$ws.Cells.Item(132, 14).Value = -67008.05929999999
$ws.Cells.Item(133, 8).Value = 6753
$ws.Cells.Item(133, 10).Value = 9999.5
$ws.Cells.Item(133, 12).Value = 29998.5
$ws.Cells.Item(133, 14).Value = -40118.5
$ws.Cells.Item(134, 8).Value = 254382.25
$ws.Cells.Item(134, 9).Value = 254382.25
$ws.Cells.Item(134, 11).Value = 763146.75
$ws.Cells.Item(134, 13).Value = -758076.75
$ws.Cells.Item(135, 8).Value = 792.4167
$ws.Cells.Item(135, 9).Value = 567.55554
$ws.Cells.Item(135, 10).Value = 1467
$ws.Cells.Item(135, 11).Value = 5107.99986
$ws.Cells.Item(135, 12).Value = 13203
$ws.Cells.Item(135, 13).Value = -2572.99986
$ws.Cells.Item(135, 14).Value = -18273
$ws.Cells.Item(137, 8).Value = 7138.769
$ws.Cells.Item(137, 10).Value = 12893
$ws.Cells.Item(137, 12).Value = 38679
$ws.Cells.Item(137, 14).Value = -48879
$ws.Cells.Item(139, 8).Value = 4786.846
$ws.Cells.Item(139, 9).Value = 2333
$ws.Cells.Item(139, 11).Value = 6999
$ws.Cells.Item(139, 13).Value = -1859
$ws.Cells.Item(140, 8).Value = 2639.8
$ws.Cells.Item(140, 9).Value = 2639.8
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 7919.400000000001
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 13).Value = -2739.400000000001
$ws.Cells.Item(140, 14).Value = $null
$ws.Cells.Item(141, 8).Value = 3113
$ws.Cells.Item(141, 9).Value = 3113
$ws.Cells.Item(141, 11).Value = 9339
$ws.Cells.Item(141, 13).Value = -4159
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 14).Value = $null
$ws.Cells.Item(113, 8).Value = 2965.5833
$ws.Cells.Item(113, 9).Value = 2740.7144
$ws.Cells.Item(113, 11).Value = 2740.7144
$ws.Cells.Item(113, 13).Value = -570.7143999999998
$ws.Cells.Item(122, 8).Value = 11338.333
$ws.Cells.Item(122, 9).Value = 14503.5
$ws.Cells.Item(122, 11).Value = 43510.5
$ws.Cells.Item(122, 13).Value = -41060.5
$ws.Cells.Item(132, 8).Value = 6265.5
$ws.Cells.Item(132, 9).Value = 6915.15
$ws.Cells.Item(132, 10).Value = 3017.25
$ws.Cells.Item(132, 11).Value = 20745.45
$ws.Cells.Item(132, 12).Value = 9051.75
$ws.Cells.Item(132, 13).Value = -18215.45
$ws.Cells.Item(132, 14).Value = -14111.75
$ws.Cells.Item(134, 8).Value = 32496
$ws.Cells.Item(134, 10).Value = 32496
$ws.Cells.Item(134, 12).Value = 97488
$ws.Cells.Item(134, 14).Value = -102558
$ws.Cells.Item(135, 8).Value = 94995.86
$ws.Cells.Item(135, 10).Value = 94995.86
$ws.Cells.Item(135, 12).Value = 94995.86
$ws.Cells.Item(135, 14).Value = -105135.86
$ws.Cells.Item(136, 8).Value = 81813.11
$ws.Cells.Item(136, 10).Value = 81813.11
$ws.Cells.Item(136, 12).Value = 245439.33
$ws.Cells.Item(136, 14).Value = -250539.33
$ws.Cells.Item(137, 8).Value = 100000
$ws.Cells.Item(137, 9).Value = 100000
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 100000
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 13).Value = -94900
$ws.Cells.Item(137, 14).Value = $null
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 3030.6667
$ws.Cells.Item(16, 9).Value = 2401.5557
$ws.Cells.Item(16, 11).Value = 2401.5557
$ws.Cells.Item(16, 13).Value = -2231.5557
$ws.Cells.Item(22, 8).Value = 130651.875
$ws.Cells.Item(22, 9).Value = 3266.0833
$ws.Cells.Item(22, 10).Value = 512809.25
$ws.Cells.Item(22, 11).Value = 3266.0833
$ws.Cells.Item(22, 12).Value = 512809.25
$ws.Cells.Item(22, 13).Value = -2971.0833
$ws.Cells.Item(22, 14).Value = -513399.25
$ws.Cells.Item(23, 8).Value = 12248.5
$ws.Cells.Item(23, 9).Value = 19997
$ws.Cells.Item(23, 10).Value = 4500
$ws.Cells.Item(23, 11).Value = 19997
$ws.Cells.Item(23, 12).Value = 4500
$ws.Cells.Item(23, 13).Value = -19767
$ws.Cells.Item(23, 14).Value = -4960
$ws.Cells.Item(27, 8).Value = 130651.875
$ws.Cells.Item(27, 9).Value = 3266.0833
$ws.Cells.Item(27, 10).Value = 512809.25
$ws.Cells.Item(27, 11).Value = 3266.0833
$ws.Cells.Item(27, 12).Value = 512809.25
$ws.Cells.Item(27, 13).Value = -3159.0833
$ws.Cells.Item(27, 14).Value = -513023.25
$ws.Cells.Item(40, 8).Value = 4632750.5
$ws.Cells.Item(40, 9).Value = 3144.5
$ws.Cells.Item(40, 11).Value = 3144.5
$ws.Cells.Item(40, 13).Value = -3008.5
$ws.Cells.Item(43, 8).Value = 33338
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 33338
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 33338
$ws.Cells.Item(43, 13).Value = $null
$ws.Cells.Item(43, 14).Value = -33724
$ws.Cells.Item(46, 8).Value = 3926.5
$ws.Cells.Item(46, 9).Value = 2999
$ws.Cells.Item(46, 10).Value = 3997.8462
$ws.Cells.Item(46, 11).Value = 2999
$ws.Cells.Item(46, 12).Value = 3997.8462
$ws.Cells.Item(46, 13).Value = -2811
$ws.Cells.Item(46, 14).Value = -4373.8462
$ws.Cells.Item(95, 8).Value = 0
$ws.Cells.Item(95, 10).Value = 0
$ws.Cells.Item(95, 12).Value = 0
$ws.Cells.Item(95, 14).Value = $null
$ws.Cells.Item(122, 8).Value = 120003590
$ws.Cells.Item(122, 9).Value = 166669980
$ws.Cells.Item(122, 10).Value = 50004004
$ws.Cells.Item(122, 11).Value = 500009940
$ws.Cells.Item(122, 12).Value = 150012012
$ws.Cells.Item(122, 13).Value = -500007490
$ws.Cells.Item(122, 14).Value = -150016912
$ws.Cells.Item(132, 8).Value = 1849.25
$ws.Cells.Item(132, 9).Value = 1799
$ws.Cells.Item(132, 10).Value = 2000
$ws.Cells.Item(132, 11).Value = 5397
$ws.Cells.Item(132, 12).Value = 6000
$ws.Cells.Item(132, 13).Value = -2867
$ws.Cells.Item(132, 14).Value = -11060
$ws.Cells.Item(136, 8).Value = 2828.4285
$ws.Cells.Item(136, 9).Value = 2622.5557
$ws.Cells.Item(136, 11).Value = 7867.6671
$ws.Cells.Item(136, 13).Value = -5317.6671
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 14).Value = $null
$ws.Cells.Item(52, 8).Value = 11009.25
$ws.Cells.Item(52, 9).Value = 10996.667
$ws.Cells.Item(52, 11).Value = 10996.667
$ws.Cells.Item(52, 13).Value = -10770.667
$ws.Cells.Item(54, 8).Value = 54395
$ws.Cells.Item(54, 10).Value = 54395
$ws.Cells.Item(54, 12).Value = 54395
$ws.Cells.Item(54, 14).Value = -55435
$ws.Cells.Item(57, 8).Value = 79000
$ws.Cells.Item(57, 9).Value = 79000
$ws.Cells.Item(57, 11).Value = 79000
$ws.Cells.Item(57, 13).Value = -78246
$ws.Cells.Item(81, 8).Value = 127575
$ws.Cells.Item(81, 9).Value = 3650
$ws.Cells.Item(81, 10).Value = 251500
$ws.Cells.Item(81, 11).Value = 7300
$ws.Cells.Item(81, 12).Value = 503000
$ws.Cells.Item(81, 13).Value = -6239
$ws.Cells.Item(81, 14).Value = -505122
$ws.Cells.Item(84, 8).Value = 127575
$ws.Cells.Item(84, 9).Value = 3650
$ws.Cells.Item(84, 10).Value = 251500
$ws.Cells.Item(84, 11).Value = 36500
$ws.Cells.Item(84, 12).Value = 2515000
$ws.Cells.Item(84, 13).Value = -31196
$ws.Cells.Item(84, 14).Value = -2525608
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 14).Value = $null
$ws.Cells.Item(96, 8).Value = 3099413.2
$ws.Cells.Item(96, 10).Value = 3762662.5
$ws.Cells.Item(96, 12).Value = 3762662.5
$ws.Cells.Item(96, 14).Value = -3765408.5
$ws.Cells.Item(100, 8).Value = 5102806.5
$ws.Cells.Item(100, 9).Value = 6494245
$ws.Cells.Item(100, 11).Value = 12988490
$ws.Cells.Item(100, 13).Value = -12987949
$ws.Cells.Item(113, 8).Value = 1233.2
$ws.Cells.Item(113, 9).Value = 1356.5714
$ws.Cells.Item(113, 11).Value = 4069.7142
$ws.Cells.Item(113, 13).Value = -1899.7142
$ws.Cells.Item(122, 8).Value = 1816.25
$ws.Cells.Item(122, 9).Value = 1816.25
$ws.Cells.Item(122, 11).Value = 5448.75
$ws.Cells.Item(122, 13).Value = -2998.75
$ws.Cells.Item(126, 8).Value = 4277.476
$ws.Cells.Item(126, 9).Value = 3517.611
$ws.Cells.Item(126, 10).Value = 8836.666999999999
$ws.Cells.Item(126, 11).Value = 10552.833
$ws.Cells.Item(126, 12).Value = 26510.001
$ws.Cells.Item(126, 13).Value = -8082.832999999999
$ws.Cells.Item(126, 14).Value = -31450.001
$ws.Cells.Item(132, 8).Value = 1950
$ws.Cells.Item(132, 9).Value = 1502
$ws.Cells.Item(132, 11).Value = 4506
$ws.Cells.Item(132, 13).Value = -1976
$ws.Cells.Item(136, 8).Value = 1758.7222
$ws.Cells.Item(136, 9).Value = 1758.7222
$ws.Cells.Item(136, 11).Value = 5276.1666
$ws.Cells.Item(136, 13).Value = -2726.1666
